$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New player roster table for rows 2-18: Name (A), Position (B), Team (C).
$data = @(
    @{Row=2;  Name="CJ McCollum";             Pos="PG,SG";    Team="New Orleans Pelicans"},
    @{Row=3;  Name="Dennis Schröder";         Pos="PG,SG";    Team="Golden State Warriors"},
    @{Row=4;  Name="Jordan Poole";            Pos="PG,SG";    Team="Washington Wizards"},
    @{Row=5;  Name="Lauri Markkanen";         Pos="SF,PF";    Team="Utah Jazz"},
    @{Row=6;  Name="Zach LaVine";             Pos="SG,SF";    Team="Chicago Bulls"},
    @{Row=7;  Name="RJ Barrett";              Pos="SG,SF,PF"; Team="Toronto Raptors"},
    @{Row=8;  Name="Nikola Jovic";            Pos="PF,C";     Team="Miami Heat"},
    @{Row=9;  Name="Kyrie Irving";            Pos="PG,SG";    Team="Dallas Mavericks"},
    @{Row=10; Name="Mike Conley";             Pos="PG";       Team="Minnesota Timberwolves"},
    @{Row=11; Name="Shai Gilgeous-Alexander"; Pos="PG,SG";    Team="Oklahoma City Thunder"},
    @{Row=12; Name="John Collins";            Pos="PF,C";     Team="Utah Jazz"},
    @{Row=13; Name="Jamal Murray";            Pos="PG,SG";    Team="Denver Nuggets"},
    @{Row=14; Name="Bam Adebayo";             Pos="C";        Team="Miami Heat"},
    @{Row=15; Name="Toumani Camara";          Pos="SF,PF";    Team="Portland Trail Blazers"},
    @{Row=16; Name="Jerami Grant";            Pos="SF,PF";    Team="Portland Trail Blazers"},
    @{Row=17; Name="Joel Embiid";             Pos="C";        Team="Philadelphia 76ers"},
    @{Row=18; Name="Tobias Harris";           Pos="SF,PF";    Team="Detroit Pistons"}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Name
    $ws.Cells.Item($r, 2).Value = $entry.Pos
    $ws.Cells.Item($r, 3).Value = $entry.Team
}
